# Weekly data refresh: a new "Acelga" price observation was inserted at
# row 360 (pushing the existing rows 360-440 down to 361-441, and
# extending the sheet's used range from A1:R440 to A1:R441).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 360; everything at/after row 360
# (including formatting) shifts down by one row.
$ws.Rows("360:360").Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(360, 1).Value2  = 3
$ws.Cells.Item(360, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(360, 3).Value2  = "Coquimbo"
$ws.Cells.Item(360, 4).Value2  = 44889
$ws.Cells.Item(360, 5).Value2  = 5
$ws.Cells.Item(360, 6).Value2  = 100112009
$ws.Cells.Item(360, 7).Value2  = "Acelga"
$ws.Cells.Item(360, 8).Value2  = "Sin especificar"
$ws.Cells.Item(360, 9).Value2  = "Primera"
$ws.Cells.Item(360, 10).Value2 = 335
$ws.Cells.Item(360, 11).Value2 = 3500
$ws.Cells.Item(360, 12).Value2 = 3800
$ws.Cells.Item(360, 13).Value2 = 3648
$ws.Cells.Item(360, 14).Value2 = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(360, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(360, 16).Value2 = 608
$ws.Cells.Item(360, 17).Value2 = 6
$ws.Cells.Item(360, 18).Value2 = "Hortaliza"
